# Update the "取得日時" (acquired timestamp) column for rows 2-6 on the
# "ランサーズ" sheet from "2025-12-31 01:26:16" to "2025-12-31 01:58:20".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-31 01:58:20"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
